# Rename the workbook's (only) sheet from "Sheet0" to "data" and move the
# active selection to D26, matching the tracked OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet "Sheet0" -> "data"
$ws.Name = "data"

# Move the active cell / selection to D26
$ws.Range("D26").Select()
